$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "['so', 'fr', 'en', 'nl', 'undetected']"
$ws.Range("I6").Value = "['tr', 'fr', 'sl', 'mk', 'en', 'id', 'ar', 'hr']"
$ws.Range("N6").Value = "['pt', 'nl', 'fr', 'de', 'uk', 'en', 'es', 'it', 'id', 'ar', 'undetected', 'ca']"
$ws.Range("R6").Value = "['es', 'ar', 'en', 'hu']"
